$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data: 98. Validate Binary Search Tree
$ws.Cells.Item(6, 1).Value = "Tree"
$ws.Cells.Item(6, 2).Value = "Medium"

$ws.Cells.Item(6, 3).Value = "98. Validate Binary Search Tree"
$ws.Cells.Item(6, 3).Style = "Good"
$ws.Cells.Item(6, 3).WrapText = $true
$ws.Cells.Item(6, 3).VerticalAlignment = -4108

$ws.Cells.Item(6, 4).Value = "Do inorder traversal and check if prev < node.val, for each node otherwise return False. For it to be a binary tree, the inorder MUST be non-increasing/non-decreasing(or ascending/descending, depends on type of bst). Checking for node.left < node.val < node.right for each node is not enough"

# Match the row height that Excel computed for the wrapped text (same as row 3)
$ws.Rows.Item(6).RowHeight = 43.2

# Hyperlink the new problem name to its Leetcode page
$url = "https://leetcode.com/problems/validate-binary-search-tree/"
$null = $ws.Hyperlinks.Add($ws.Range("C6"), $url, "", "", $url)

# Update selection to reflect where the user ended up after editing
$null = $ws.Range("D8").Select()
